$wb = $excel.ActiveWorkbook

# --- Rename the "User Name" header cell on the existing sheet to "UserName" ---
$ws1 = $wb.Worksheets.Item("invalidCredentialTest")
$ws1.Range("A1").Value = "UserName"

# --- Add two new, empty worksheets at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet1 = $wb.Worksheets.Add($null, $lastSheet)
$sheet1.Name = "Sheet1"
$sheet2 = $wb.Worksheets.Add($null, $sheet1)
$sheet2.Name = "Sheet2"

# --- Re-activate the original sheet and select the entire row 5 ---
$ws1.Activate()
$ws1.Range("A5:XFD1048576").Select()
